# Applies the commit "remove 'timechop' + capitalize 'AS'" to the deck.
#
# Three textual edits, one per slide:
#   Slide 10 (SQL code block):  "... THEN 1 ELSE 0 END as label" -> "... END AS label"
#   Slide 35 (title):           "Configuring Temporal Parameters (Timechop)" -> "Configuring Temporal Parameters"
#   Slide 43 (bullet body):     merge two runs of the "Train and test ..." bullet into one run

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 10: capitalize "as" -> "AS" in "... THEN 1 ELSE 0 END as label"
# ---------------------------------------------------------------------------
$slide10 = $p.Slides.Item(10)
for ($i = 1; $i -le $slide10.Shapes.Count; $i++) {
    $shape = $slide10.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        $full = $tr.Text
        if ($full.Contains("THEN 1 ELSE 0 END as label")) {
            $idx = $full.IndexOf("as label")
            # Select "as " (including trailing space) and retype as "AS "
            $sub = $tr.Characters($idx + 1, 3)
            $sub.Text = "AS "
        }
    }
}

# ---------------------------------------------------------------------------
# Slide 35: drop " (Timechop)" from the title
# ---------------------------------------------------------------------------
$slide35 = $p.Slides.Item(35)
for ($i = 1; $i -le $slide35.Shapes.Count; $i++) {
    $shape = $slide35.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        $full = $tr.Text
        if ($full.Contains("Configuring Temporal Parameters (Timechop)")) {
            # Retype "Parameters" in place (splits the run at the right boundary)
            $idxParam = $full.IndexOf("Parameters")
            $lenParam = "Parameters".Length
            $subParam = $tr.Characters($idxParam + 1, $lenParam)
            $subParam.Text = "Parameters"

            # Remove the trailing " (Timechop)"
            $full2 = $tr.Text
            $idxSuffix = $full2.IndexOf(" (Timechop)")
            $lenSuffix = " (Timechop)".Length
            $subSuffix = $tr.Characters($idxSuffix + 1, $lenSuffix)
            $subSuffix.Text = ""
        }
    }
}

# ---------------------------------------------------------------------------
# Slide 43: merge "Train and test labels aggregate data from " + "overlapping
# times" into a single run
# ---------------------------------------------------------------------------
$slide43 = $p.Slides.Item(43)
for ($i = 1; $i -le $slide43.Shapes.Count; $i++) {
    $shape = $slide43.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        $full = $tr.Text
        $needle = "Train and test labels aggregate data from overlapping times"
        if ($full.Contains($needle)) {
            $idx = $full.IndexOf($needle)
            $len = $needle.Length
            $sub = $tr.Characters($idx + 1, $len)
            $sub.Text = $needle
        }
    }
}
